$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Price (column D) / Volume(1h) (column E) readings for this
# run of the cryptos scraper. Map of cell address -> new text value.
$newValues = [ordered]@{
    'D2' = '29.334.71'
    'E2' = '  -0.44%  '
    'D3' = '1.846.09'
    'E3' = '  -0.34%  '
    'D4' = '0.9977'
    'E4' = '  -0.31%  '
    'D5' = '240.16'
    'E5' = '  -0.30%  '
    'D6' = '0.6283'
    'E6' = '  -0.37%  '
    'D7' = '0.9997'
    'E7' = '  -0.15%  '
    'D8' = '0.07571'
    'E8' = '  -1.17%  '
    'D9' = '0.2909'
    'E9' = '  -1.05%  '
    'D10' = '24.59'
    'E10' = '  -0.23%  '
    'D11' = '0.07736'
    'E11' = '  -0.27%  '
    'D12' = '1.847.05'
    'E12' = '  -0.60%  '
    'D13' = '5.013'
    'E13' = '  -0.39%  '
    'D14' = '0.6787'
    'E14' = '  -0.30%  '
    'E15' = '  -3.02%  '
    'D16' = '82.99'
    'E16' = '  -0.77%  '
    'D17' = '6.102'
    'E17' = '  -1.11%  '
    'D18' = '29.323.82'
    'E18' = '  -0.52%  '
    'D19' = '228.89'
    'E19' = '  -0.15%  '
    'D20' = '12.33'
    'E20' = '  -1.15%  '
    'D21' = '0.9993'
    'E21' = '  -0.16%  '
    'D22' = '7.423'
    'E22' = '  -0.38%  '
    'D23' = '1.000'
    'E23' = '  -0.15%  '
    'D24' = '158.82'
    'E24' = '  +1.11%  '
    'D25' = '0.1387'
    'E25' = '  +0.14%  '
    'D26' = '8.425'
    'D27' = '17.62'
    'E27' = '  -0.56%  '
    'D28' = '1.430'
    'E28' = '  +7.87%  '
    'D30' = '0.05662'
    'E30' = '  -0.38%  '
    'D31' = '4.105'
    'E31' = '  -0.69%  '
    'D32' = '4.034'
    'E32' = '  -0.32%  '
    'D33' = '1.154'
    'E33' = '  -1.01%  '
    'D34' = '1.818'
    'E34' = '  -1.73%  '
    'D35' = '0.6945'
    'E35' = '  -1.07%  '
    'D36' = '2.579'
    'E36' = '  -0.30%  '
    'D37' = '0.01830'
    'D38' = '1.234.90'
    'E38' = '  +1.36%  '
    'D39' = '2.719'
    'E39' = '  -2.36%  '
    'D40' = '6.373'
    'E40' = '  -2.75%  '
    'D41' = '0.8951'
    'E41' = '  -1.52%  '
    'D42' = '0.9988'
    'E42' = '  -0.29%  '
    'D43' = '101.21'
    'E43' = '  -0.47%  '
    'D44' = '65.36'
    'E44' = '  -1.63%  '
    'E45' = '  -1.23%  '
    'D46' = '7.130'
    'E46' = '  +0.32%  '
    'D47' = '0.3996'
    'E47' = '  -0.60%  '
    'D48' = '0.1150'
    'E48' = '  +1.06%  '
    'D49' = '8.956'
    'E49' = '  -0.55%  '
    'D50' = '1.674'
    'E50' = '  -0.59%  '
    'D51' = '0.05693'
    'E51' = '  -0.46%  '
}

foreach ($addr in $newValues.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking values (e.g. "240.16") stay
    # text instead of being auto-converted to numbers, matching the
    # source cell type; restoring the default style afterwards keeps the
    # cell formatting identical to the original (unstyled) cells.
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$addr]
    $cell.Style = "Normal"
}
